$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 517, shifting existing rows (517..578) down to (518..579).
$ws.Rows(517).Insert()

# Populate the newly inserted row 517 with the new record's data.
$ws.Cells.Item(517, 1).Value  = 6
$ws.Cells.Item(517, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(517, 3).Value  = "Metropolitana"
$ws.Cells.Item(517, 4).Value  = 45077
$ws.Cells.Item(517, 5).Value  = 13
$ws.Cells.Item(517, 6).Value  = 100112032
$ws.Cells.Item(517, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(517, 8).Value  = "Sin especificar"
$ws.Cells.Item(517, 9).Value  = "Primera"
$ws.Cells.Item(517, 10).Value = 800
$ws.Cells.Item(517, 11).Value = 7000
$ws.Cells.Item(517, 12).Value = 8000
$ws.Cells.Item(517, 13).Value = 7625
$ws.Cells.Item(517, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(517, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(517, 16).Value = 152
$ws.Cells.Item(517, 17).Value = 50
$ws.Cells.Item(517, 18).Value = "Hortaliza"
